$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D from Excel's automatic text->number coercion
# (values like "1.000" / "49.30" would otherwise become numeric).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.425.71"
$ws.Range("E2").Value = "  -3.87%  "

$ws.Range("D3").Value = "1.572.75"
$ws.Range("E3").Value = "  -3.44%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "289.90"
$ws.Range("E6").Value = "  -2.64%  "

$ws.Range("D7").Value = "0.3673"
$ws.Range("E7").Value = "  -2.43%  "

$ws.Range("D8").Value = "49.30"
$ws.Range("E8").Value = "  -1.68%  "

$ws.Range("D9").Value = "0.3404"
$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("D10").Value = "1.173"
$ws.Range("E10").Value = "  -2.23%  "

$ws.Range("D11").Value = "0.07650"
$ws.Range("E11").Value = "  -4.71%  "

$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").Value = "21.28"
$ws.Range("E13").Value = "  -2.69%  "

$ws.Range("D14").Value = "6.073"
$ws.Range("E14").Value = "  -3.42%  "

$ws.Range("D15").Value = "6.927"
$ws.Range("E15").Value = "  -4.03%  "

$ws.Range("D16").Value = "1.572.39"
$ws.Range("E16").Value = "  -3.15%  "

$ws.Range("E17").Value = "  -4.48%  "

$ws.Range("D18").Value = "89.95"
$ws.Range("E18").Value = "  -4.86%  "

$ws.Range("D19").Value = "0.06743"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").Value = "6.263"
$ws.Range("E21").Value = "  -5.05%  "

$ws.Range("D22").Value = "16.60"
$ws.Range("E22").Value = "  -3.67%  "

$ws.Range("E23").Value = "  -5.78%  "

$ws.Range("D24").Value = "12.02"
$ws.Range("E24").Value = "  -2.52%  "

$ws.Range("D25").Value = "22.415.08"
$ws.Range("E25").Value = "  -3.92%  "

$ws.Range("D26").Value = "2.387"
$ws.Range("E26").Value = "  -1.71%  "

$ws.Range("D27").Value = "2.903"
$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  -4.04%  "

$ws.Range("D29").Value = "146.18"
$ws.Range("E29").Value = "  -2.95%  "

$ws.Range("D30").Value = "4.989"
$ws.Range("E30").Value = "  -3.30%  "

$ws.Range("D31").Value = "125.73"
$ws.Range("E31").Value = "  -4.44%  "

$ws.Range("D32").Value = "1.741.66"
$ws.Range("E32").Value = "  -3.49%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.032"
$ws.Range("E33").Value = "  +5.93%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "6.252"
$ws.Range("E34").Value = "  -6.71%  "

$ws.Range("D35").Value = "2.020"
$ws.Range("E35").Value = "  -4.53%  "

$ws.Range("D36").Value = "10.15"
$ws.Range("E36").Value = "  -9.06%  "

$ws.Range("D37").Value = "0.08465"
$ws.Range("E37").Value = "  -3.15%  "

$ws.Range("D38").Value = "0.02545"
$ws.Range("E38").Value = "  -4.10%  "

$ws.Range("D39").Value = "0.2326"
$ws.Range("E39").Value = "  -3.68%  "

$ws.Range("D40").Value = "5.535"
$ws.Range("E40").Value = "  -5.08%  "

$ws.Range("D41").Value = "0.06499"
$ws.Range("E41").Value = "  -2.74%  "

$ws.Range("D42").Value = "1.298"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("E43").Value = "  -7.20%  "

$ws.Range("D44").Value = "0.6373"
$ws.Range("E44").Value = "  -6.01%  "

$ws.Range("D45").Value = "14.27"
$ws.Range("E45").Value = "  -6.93%  "

$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").Value = "0.6007"
$ws.Range("E47").Value = "  -4.60%  "

$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("D49").Value = "2.114"
$ws.Range("E49").Value = "  -5.10%  "

$ws.Range("D50").Value = "1.256"
$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("D51").Value = "124.79"
$ws.Range("E51").Value = "  -0.98%  "

# Restore default (General) style now that text values are safely stored,
# without leaving the explicit Text number format applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
